# Update crypto price/volume data per the Aug 30 2023 GitHub Actions refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper (scratch) cell used to push pure-numeric-looking strings into cells
# as literal text (so "224.29" stays the text "224.29", not the float 224.29).
# We build it as a formula returning a text string, copy it, and paste-special
# values-only into the destination so the destination cell keeps its original
# (default) style -- only the stored value changes, exactly like the source diff.
$helper = $ws.Range("ZZ1")

$ws.Range("D2").Value = "27.379.52"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.710.49"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.15%  "
$helper.Formula = '="224.29"'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.71%  "
$helper.Formula = '="0.5338"'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -2.47%  "
$helper.Formula = '="0.2666"'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -4.12%  "
$helper.Formula = '="0.06598"'
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -2.36%  "
$helper.Formula = '="20.88"'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -4.51%  "
$helper.Formula = '="0.07645"'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -1.83%  "
$helper.Formula = '="4.572"'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "1.728.99"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "1.946.37"
$helper.Formula = '="0.5757"'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").Value = "0.0₅8173"
$ws.Range("E16").Value = "  -2.95%  "
$helper.Formula = '="67.85"'
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "27.359.14"
$ws.Range("E18").Value = "  -2.09%  "
$helper.Formula = '="215.67"'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("E20").Value = "  +0.13%  "
$helper.Formula = '="4.675"'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  -4.34%  "
$helper.Formula = '="5.977"'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -4.47%  "
$helper.Formula = '="1.004"'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.09%  "
$helper.Formula = '="142.14"'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -2.70%  "
$helper.Formula = '="1.721"'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +2.08%  "
$helper.Formula = '="0.1217"'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.99%  "
$helper.Formula = '="7.280"'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -2.73%  "
$helper.Formula = '="16.32"'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -5.31%  "
$helper.Formula = '="0.05410"'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -4.96%  "
$helper.Formula = '="1.294"'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -1.77%  "
$helper.Formula = '="3.506"'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -5.48%  "
$helper.Formula = '="3.432"'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("E34").Value = "  -2.73%  "
$helper.Formula = '="2.881"'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.80%  "
$helper.Formula = '="0.9495"'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("E37").Value = "  -1.23%  "
$helper.Formula = '="0.5868"'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.14%  "
$helper.Formula = '="0.01632"'
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$helper.Formula = '="5.868"'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "1.044.21"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  +0.10%  "
$helper.Formula = '="0.8416"'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.20%  "
$helper.Formula = '="100.87"'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "1.852.98"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  +1.58%  "
$helper.Formula = '="58.10"'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -2.89%  "
$helper.Formula = '="0.4510"'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +1.56%  "
$helper.Formula = '="1.004"'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.07%  "
$helper.Formula = '="8.082"'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -2.54%  "
$helper.Formula = '="0.05241"'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.78%  "

# Remove the scratch cell entirely (shift-up delete leaves no trace behind).
$helper.Delete(-4159)
$excel.CutCopyMode = $false
